$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1176.1818
$ws.Range("I19").Value = 1093.875
$ws.Range("K19").Value = 1093.875
$ws.Range("M19").Value = -918.875
$ws.Range("H41").Value = 799.4666999999999
$ws.Range("I41").Value = 1148.3
$ws.Range("K41").Value = 1148.3
$ws.Range("M41").Value = -708.3
$ws.Range("H58").Value = 4802.4287
$ws.Range("J58").Value = 7749.25
$ws.Range("L58").Value = 23247.75
$ws.Range("N58").Value = -23547.75
$ws.Range("H74").Value = 7960.75
$ws.Range("I74").Value = 7960.75
$ws.Range("K74").Value = 7960.75
$ws.Range("M74").Value = -7024.75
$ws.Range("H77").Value = 7960.75
$ws.Range("I77").Value = 7960.75
$ws.Range("K77").Value = 39803.75
$ws.Range("M77").Value = -35123.75
$ws.Range("H80").Value = 2007.7
$ws.Range("I80").Value = 3027.6667
$ws.Range("K80").Value = 9083.000100000001
$ws.Range("M80").Value = -8085.000100000001
$ws.Range("H83").Value = 2007.7
$ws.Range("I83").Value = 3027.6667
$ws.Range("K83").Value = 27249.0003
$ws.Range("M83").Value = -22257.0003
$ws.Range("H86").Value = 9764.714
$ws.Range("I86").Value = 11687.25
$ws.Range("J86").Value = 7201.3335
$ws.Range("K86").Value = 11687.25
$ws.Range("L86").Value = 7201.3335
$ws.Range("M86").Value = -10564.25
$ws.Range("N86").Value = -9447.333500000001
$ws.Range("H89").Value = 9764.714
$ws.Range("I89").Value = 11687.25
$ws.Range("J89").Value = 7201.3335
$ws.Range("K89").Value = 58436.25
$ws.Range("L89").Value = 36006.6675
$ws.Range("M89").Value = -52820.25
$ws.Range("N89").Value = -47238.6675
$ws.Range("H97").Value = 3769.8462
$ws.Range("I97").Value = 903
$ws.Range("J97").Value = 4008.75
$ws.Range("K97").Value = 2709
$ws.Range("L97").Value = 12026.25
$ws.Range("M97").Value = -2213
$ws.Range("N97").Value = -13018.25
$ws.Range("H107").Value = 1450
$ws.Range("I107").Value = 1450
$ws.Range("K107").Value = 1450
$ws.Range("M107").Value = 470
$ws.Range("H116").Value = 7500
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 7500
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 7500
$ws.Range("M116").ClearContents() | Out-Null
$ws.Range("N116").Value = -14384
$ws.Range("H129").Value = 14484.714
$ws.Range("I129").Value = 2079.2
$ws.Range("K129").Value = 6237.599999999999
$ws.Range("M129").Value = -1237.599999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5676
$ws.Range("I45").Value = 8557.111000000001
$ws.Range("K45").Value = 8557.111000000001
$ws.Range("M45").Value = -8180.111000000001
$ws.Range("H61").Value = 76925270
$ws.Range("I61").Value = 100002100
$ws.Range("K61").Value = 100002100
$ws.Range("M61").Value = -100001888
$ws.Range("H110").Value = 53783.367
$ws.Range("I110").Value = 67900.336
$ws.Range("J110").Value = 844.75
$ws.Range("K110").Value = 67900.336
$ws.Range("L110").Value = 844.75
$ws.Range("M110").Value = -65855.336
$ws.Range("N110").Value = -4934.75
$ws.Range("H122").Value = 7481.9165
$ws.Range("I122").Value = 6848.5
$ws.Range("K122").Value = 20545.5
$ws.Range("M122").Value = -18095.5
$ws.Range("H136").Value = 76925270
$ws.Range("I136").Value = 100002100
$ws.Range("K136").Value = 300006300
$ws.Range("M136").Value = -300003750
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 40143.4
$ws.Range("J6").Value = 46474.332
$ws.Range("L6").Value = 46474.332
$ws.Range("N6").Value = -46700.332
$ws.Range("H94").Value = 2524.4443
$ws.Range("I94").Value = 2359.4546
$ws.Range("K94").Value = 2359.4546
$ws.Range("M94").Value = -1908.4546
$ws.Range("H107").Value = 88795.836
$ws.Range("I107").Value = 5395.6665
$ws.Range("K107").Value = 5395.6665
$ws.Range("M107").Value = -3475.6665
$ws.Range("H109").Value = 66666
$ws.Range("J109").Value = 66666
$ws.Range("L109").Value = 66666
$ws.Range("N109").Value = -69440
$ws.Range("H112").Value = 56500
$ws.Range("J112").Value = 56500
$ws.Range("L112").Value = 56500
$ws.Range("N112").Value = -59454
$ws.Range("H117").Value = 26998
$ws.Range("J117").Value = 26998
$ws.Range("L117").Value = 26998
$ws.Range("N117").Value = -36176
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11854.3
$ws.Range("J86").Value = 14108.728
$ws.Range("L86").Value = 14108.728
$ws.Range("N86").Value = -16354.728
$ws.Range("H89").Value = 11854.3
$ws.Range("J89").Value = 14108.728
$ws.Range("L89").Value = 70543.64
$ws.Range("N89").Value = -81775.64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 48
$ws.Range("I40").Value = 48
$ws.Range("K40").Value = 192
$ws.Range("M40").Value = -123
$ws.Range("H68").Value = 2875
$ws.Range("J68").Value = 10000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31622
$ws.Range("H71").Value = 2875
$ws.Range("J71").Value = 10000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -98112
$ws.Range("H113").Value = 111772.22
$ws.Range("J113").Value = 578.4286
$ws.Range("L113").Value = 1735.2858
$ws.Range("N113").Value = -6075.2858
$ws.Range("H131").Value = 1710.6471
$ws.Range("I131").Value = 1138.1875
$ws.Range("K131").Value = 3414.5625
$ws.Range("M131").Value = 1625.4375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 78999
$ws.Range("J120").Value = 78999
$ws.Range("L120").Value = 78999
$ws.Range("N120").Value = -88675
$ws.Range("H122").Value = 5104.7
$ws.Range("I122").Value = 3116.3333
$ws.Range("K122").Value = 9348.999899999999
$ws.Range("M122").Value = -6898.999899999999
$ws.Range("H132").Value = 25002400
$ws.Range("I132").Value = 31252250
$ws.Range("K132").Value = 93756750
$ws.Range("M132").Value = -93754220
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2242.5
$ws.Range("I16").Value = 1527.091
$ws.Range("J16").Value = 4865.6665
$ws.Range("K16").Value = 1527.091
$ws.Range("L16").Value = 4865.6665
$ws.Range("M16").Value = -1357.091
$ws.Range("N16").Value = -5205.6665
$ws.Range("H40").Value = 3748.1667
$ws.Range("I40").Value = 3748.1667
$ws.Range("K40").Value = 3748.1667
$ws.Range("M40").Value = -3612.1667
$ws.Range("H46").Value = 1936.9166
$ws.Range("I46").Value = 1794.8182
$ws.Range("K46").Value = 1794.8182
$ws.Range("M46").Value = -1606.8182
$ws.Range("H61").Value = 6421.2
$ws.Range("I61").Value = 6501.4443
$ws.Range("K61").Value = 6501.4443
$ws.Range("M61").Value = -6299.4443
$ws.Range("H113").Value = 6421.2
$ws.Range("I113").Value = 6501.4443
$ws.Range("K113").Value = 6501.4443
$ws.Range("M113").Value = -4331.4443
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1459.4667
$ws.Range("I96").Value = 1114.625
$ws.Range("J96").Value = 1853.5714
$ws.Range("K96").Value = 1114.625
$ws.Range("L96").Value = 1853.5714
$ws.Range("M96").Value = 258.375
$ws.Range("N96").Value = -4599.5714
$ws.Range("H122").Value = 2272.25
$ws.Range("I122").Value = 2199.5
$ws.Range("K122").Value = 6598.5
$ws.Range("M122").Value = -4148.5
